$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in numeric budget data for Jan (B), Feb (C), Mar (D) columns
# Row 4: Rent
$ws.Range("B4").Value = 1200
$ws.Range("C4").Value = 1200
$ws.Range("D4").Value = 1000

# Row 5: Phone
$ws.Range("B5").Value = 100
$ws.Range("C5").Value = 135
$ws.Range("D5").Value = 100

# Row 6: Credit Cards
$ws.Range("B6").Value = 150
$ws.Range("C6").Value = 200
$ws.Range("D6").Value = 125

# Row 7: Food
$ws.Range("B7").Value = 300
$ws.Range("C7").Value = 275
$ws.Range("D7").Value = 350

# Row 8: Candy
$ws.Range("B8").Value = 100
$ws.Range("C8").Value = 100
$ws.Range("D8").Value = 110

# Update the active selection to E5
$ws.Range("E5").Select()
